{"js": "const replacements = [\n  { oldText: \"2024-08-06 Tuesday\", newText: \"2024-08-07 Wednesday\" },\n  { oldText: \"715\u00d76=4290\", newText: \"486\u00d77=3402\" },\n  { oldText: \"138\u00d79=1242\", newText: \"244\u00d72=488\" },\n  { oldText: \"667\u00d72=1334\", newText: \"559\u00d75=2795\" },\n  { oldText: \"713\u00d79=6417\", newText: \"116\u00d74=464\" },\n  { oldText: \"662\u00d73=1986\", newText: \"533\u00d76=3198\" },\n  { oldText: \"676\u00d75=3380\", newText: \"820\u00d72=1640\" },\n  { oldText: \"206\u00d74=824\", newText: \"362\u00d73=1086\" },\n  { oldText: \"807\u00d75=4035\", newText: \"329\u00d77=2303\" },\n  { oldText: \"529\u00d74=2116\", newText: \"633\u00d73=1899\" },\n  { oldText: \"259\u00d73=777\", newText: \"161\u00d76=966\" },\n  { oldText: \"104\u00d74=416\", newText: \"157\u00d78=1256\" },\n  { oldText: \"194\u00d79=1746\", newText: \"647\u00d73=1941\" },\n  { oldText: \"359\u00d73=1077\", newText: \"632\u00d77=4424\" },\n  { oldText: \"368\u00d75=1840\", newText: \"979\u00d77=6853\" },\n  { oldText: \"961\u00d75=4805\", newText: \"365\u00d74=1460\" },\n  { oldText: \"115\u00d76=690\", newText: \"977\u00d72=1954\" },\n  { oldText: \"587\u00d74=2348\", newText: \"917\u00d75=4585\" },\n  { oldText: \"176\u00d76=1056\", newText: \"446\u00d76=2676\" },\n  { oldText: \"783\u00d79=7047\", newText: \"786\u00d77=5502\" },\n  { oldText: \"448\u00d73=1344\", newText: \"199\u00d77=1393\" },\n  { oldText: \"625\u00d72=1250\", newText: \"374\u00d72=748\" },\n  { oldText: \"839\u00d77=5873\", newText: \"922\u00d76=5532\" },\n  { oldText: \"534\u00d75=2670\", newText: \"679\u00d77=4753\" },\n  { oldText: \"981\u00d77=6867\", newText: \"378\u00d79=3402\" },\n  { oldText: \"469\u00d79=4221\", newText: \"639\u00d78=5112\" },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"2024-08-06 Tuesday\"; new = \"2024-08-07 Wednesday\" },\n    @{ old = \"715\u00d76=4290\"; new = \"486\u00d77=3402\" },\n    @{ old = \"138\u00d79=1242\"; new = \"244\u00d72=488\" },\n    @{ old = \"667\u00d72=1334\"; new = \"559\u00d75=2795\" },\n    @{ old = \"713\u00d79=6417\"; new = \"116\u00d74=464\" },\n    @{ old = \"662\u00d73=1986\"; new = \"533\u00d76=3198\" },\n    @{ old = \"676\u00d75=3380\"; new = \"820\u00d72=1640\" },\n    @{ old = \"206\u00d74=824\"; new = \"362\u00d73=1086\" },\n    @{ old = \"807\u00d75=4035\"; new = \"329\u00d77=2303\" },\n    @{ old = \"529\u00d74=2116\"; new = \"633\u00d73=1899\" },\n    @{ old = \"259\u00d73=777\"; new = \"161\u00d76=966\" },\n    @{ old = \"104\u00d74=416\"; new = \"157\u00d78=1256\" },\n    @{ old = \"194\u00d79=1746\"; new = \"647\u00d73=1941\" },\n    @{ old = \"359\u00d73=1077\"; new = \"632\u00d77=4424\" },\n    @{ old = \"368\u00d75=1840\"; new = \"979\u00d77=6853\" },\n    @{ old = \"961\u00d75=4805\"; new = \"365\u00d74=1460\" },\n    @{ old = \"115\u00d76=690\"; new = \"977\u00d72=1954\" },\n    @{ old = \"587\u00d74=2348\"; new = \"917\u00d75=4585\" },\n    @{ old = \"176\u00d76=1056\"; new = \"446\u00d76=2676\" },\n    @{ old = \"783\u00d79=7047\"; new = \"786\u00d77=5502\" },\n    @{ old = \"448\u00d73=1344\"; new = \"199\u00d77=1393\" },\n    @{ old = \"625\u00d72=1250\"; new = \"374\u00d72=748\" },\n    @{ old = \"839\u00d77=5873\"; new = \"922\u00d76=5532\" },\n    @{ old = \"534\u00d75=2670\"; new = \"679\u00d77=4753\" },\n    @{ old = \"981\u00d77=6867\"; new = \"378\u00d79=3402\" },\n    @{ old = \"469\u00d79=4221\"; new = \"639\u00d78=5112\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"No match found for: $($r.old)\"\n    }\n}\n"}
